$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 39995
$ws.Cells.Item(2, 2).Value = "Heloísa Borges"
$ws.Cells.Item(2, 3).Value = "Financeiro"
$ws.Cells.Item(2, 5).Value = 8
$ws.Cells.Item(2, 6).Value = (Get-Date -Year 2023 -Month 6 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(2, 7).Value = 6124.57

# Row 3
$ws.Cells.Item(3, 1).Value = 88252
$ws.Cells.Item(3, 2).Value = "Dr. Erick Rodrigues"
$ws.Cells.Item(3, 3).Value = "Financeiro"
$ws.Cells.Item(3, 5).Value = 8
$ws.Cells.Item(3, 6).Value = (Get-Date -Year 2023 -Month 6 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(3, 7).Value = 8709.98

# Row 4
$ws.Cells.Item(4, 1).Value = 29869
$ws.Cells.Item(4, 2).Value = "Sara Gonçalves"
$ws.Cells.Item(4, 3).Value = "P&D"
$ws.Cells.Item(4, 4).Value = "Doenca"
$ws.Cells.Item(4, 5).Value = 6
$ws.Cells.Item(4, 6).Value = (Get-Date -Year 2023 -Month 6 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(4, 7).Value = 6934.36

# Row 5
$ws.Cells.Item(5, 1).Value = 39310
$ws.Cells.Item(5, 2).Value = "Ana Vitória Fogaça"
$ws.Cells.Item(5, 3).Value = "Recursos Humanos"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = (Get-Date -Year 2023 -Month 6 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(5, 7).Value = 5026.57

# Row 6
$ws.Cells.Item(6, 1).Value = 24118
$ws.Cells.Item(6, 2).Value = "Luan Dias"
$ws.Cells.Item(6, 3).Value = "Juridico"
$ws.Cells.Item(6, 4).Value = "Outros"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = (Get-Date -Year 2023 -Month 6 -Day 24 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(6, 7).Value = 5820.62

# Row 7
$ws.Cells.Item(7, 1).Value = 97818
$ws.Cells.Item(7, 2).Value = "Otto da Cruz"
$ws.Cells.Item(7, 3).Value = "Engenharia"
$ws.Cells.Item(7, 4).Value = "Viagem de negocios"
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = (Get-Date -Year 2023 -Month 6 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(7, 7).Value = 3591.59

# Row 8
$ws.Cells.Item(8, 1).Value = 48574
$ws.Cells.Item(8, 2).Value = "Anna Liz Casa Grande"
$ws.Cells.Item(8, 3).Value = "Engenharia"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = (Get-Date -Year 2023 -Month 6 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(8, 7).Value = 8122.26

# Row 9
$ws.Cells.Item(9, 1).Value = 60300
$ws.Cells.Item(9, 2).Value = "Sra. Isadora Sá"
$ws.Cells.Item(9, 3).Value = "TI"
$ws.Cells.Item(9, 4).Value = "Viagem de negocios"
$ws.Cells.Item(9, 5).Value = 4
$ws.Cells.Item(9, 6).Value = (Get-Date -Year 2023 -Month 6 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(9, 7).Value = 7134.06

# Row 10
$ws.Cells.Item(10, 1).Value = 35936
$ws.Cells.Item(10, 2).Value = "Dr. Isaac Freitas"
$ws.Cells.Item(10, 3).Value = "Operacoes"
$ws.Cells.Item(10, 5).Value = 8
$ws.Cells.Item(10, 6).Value = (Get-Date -Year 2023 -Month 6 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(10, 7).Value = 8506.61

# Row 11
$ws.Cells.Item(11, 1).Value = 31740
$ws.Cells.Item(11, 2).Value = "Luara Barros"
$ws.Cells.Item(11, 3).Value = "Juridico"
$ws.Cells.Item(11, 4).Value = "Doenca"
$ws.Cells.Item(11, 5).Value = 7
$ws.Cells.Item(11, 6).Value = (Get-Date -Year 2023 -Month 6 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(11, 7).Value = 4279.83
